$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Cilantro (Femacal de La Calera).
# It belongs right above the existing row 312, so push that row (and every
# row below it) down by one before filling in the new data.
$ws.Rows("312:312").Insert()

$ws.Range("A312").Value = 3
$ws.Range("B312").Value = "Femacal de La Calera"
$ws.Range("C312").Value = "Coquimbo"
$ws.Range("D312").Value = 44809
$ws.Range("E312").Value = 5
$ws.Range("F312").Value = 100112040
$ws.Range("G312").Value = "Cilantro"
$ws.Range("H312").Value = "Sin especificar"
$ws.Range("I312").Value = "Primera"
$ws.Range("J312").Value = 105
$ws.Range("K312").Value = 4500
$ws.Range("L312").Value = 4500
$ws.Range("M312").Value = 4500
$ws.Range("N312").Value = "$/docena de atados (3 kilos)"
$ws.Range("O312").Value = "Provincia de Quillota"
$ws.Range("P312").Value = 1500
$ws.Range("Q312").Value = 3
$ws.Range("R312").Value = "Hortaliza"
